# Thu Jan  5 17:57:05 UTC 2023 symbol-list refresh (GitHub Actions).
#
# Source data is scraped crypto prices/volumes; the sheet stores every
# value as literal text (e.g. "256.62", "-1.27%") rather than numbers, so
# each cell below is written with a leading apostrophe. That forces Excel
# to keep the quote-prefixed text exactly as typed instead of silently
# re-interpreting number- or percent-looking strings as numeric values.
# Plain, non-numeric text (coin names, URLs) is written as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = "'256.47"
$ws.Range("E2").Value = "'-1.29%"

# Row 3 - OKB
$ws.Range("D3").Value = "'27.23"
$ws.Range("E3").Value = "'-2.78%"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "'4.559"
$ws.Range("E4").Value = "'-12.60%"

# Row 5 - Cronos
$ws.Range("D5").Value = "'0.05894"
$ws.Range("E5").Value = "'-0.75%"

# Row 6 - KuCoinToken
$ws.Range("D6").Value = "'6.629"
$ws.Range("E6").Value = "'-1.52%"

# Row 7 - MXToken
$ws.Range("D7").Value = "'0.8586"
$ws.Range("E7").Value = "'-1.69%"

# Row 8 - FTXToken
$ws.Range("D8").Value = "'0.9265"
$ws.Range("E8").Value = "'-8.01%"

# Row 9 - WazirX
$ws.Range("D9").Value = "'0.1408"
$ws.Range("E9").Value = "'-1.35%"

# Row 10 - LiechtensteinCryptoassetsExchange
$ws.Range("D10").Value = "'0.03627"
$ws.Range("E10").Value = "'-0.42%"

# Row 11 - MandalaExchangeToken
$ws.Range("D11").Value = "'0.07099"
$ws.Range("E11").Value = "'-2.04%"

# Row 12 - BitrueCoin
$ws.Range("D12").Value = "'0.03234"
$ws.Range("E12").Value = "'1.01%"

# Row 13 - BitMartToken
$ws.Range("D13").Value = "'0.09210"
$ws.Range("E13").Value = "'-0.42%"

# Row 14 - BitForexToken
$ws.Range("D14").Value = "'0.001544"
$ws.Range("E14").Value = "'-0.02%"

# Row 15 - was TigerCash, now One (ranking shuffle)
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006061"
$ws.Range("E15").Value = "'-0.18%"

# Row 16 - was LEO, now TigerCash
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006091"
$ws.Range("E16").Value = "'3.61%"

# Row 17 - was GateToken, now LEO
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.515"
$ws.Range("E17").Value = "'0.70%"

# Row 18 - was BTSEToken, now GateToken
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.189"
$ws.Range("E18").Value = "'-1.31%"

# Row 19 - was One, now BTSEToken
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.202"
$ws.Range("E19").Value = "'-0.23%"

# Row 20 - BitpandaEcosystemToken (volume only)
$ws.Range("E20").Value = "'-2.16%"

# Row 21 - ProBitToken (volume only)
$ws.Range("E21").Value = "'-1.00%"

# Row 22 - MCDex
$ws.Range("D22").Value = "'3.848"
$ws.Range("E22").Value = "'9.16%"

# Row 23 - CoinExToken
$ws.Range("D23").Value = "'0.04211"
$ws.Range("E23").Value = "'0.74%"

# Row 24 - BitKan
$ws.Range("D24").Value = "'0.001221"
$ws.Range("E24").Value = "'0.43%"

# Row 25 - HotbitToken
$ws.Range("D25").Value = "'0.004278"
$ws.Range("E25").Value = "'-6.40%"

# Row 26 - NitroEx (volume only)
$ws.Range("E26").Value = "'0.18%"

# Row 27 - UpBots (volume only)
$ws.Range("E27").Value = "'-21.94%"

# Row 40 - IDEX
$ws.Range("D40").Value = "'0.03833"
$ws.Range("E40").Value = "'-0.45%"

# Row 41 - was KickToken, now BKEXToken
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1100"
$ws.Range("E41").Value = "'-0.96%"

# Row 42 - was BKEXToken, now KickToken
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.003950"
$ws.Range("E42").Value = "'-27.03%"

# Row 43 - CEJI
$ws.Range("D43").Value = "'0.002421"
$ws.Range("E43").Value = "'2.00%"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "'0.01137"
$ws.Range("E44").Value = "'4.30%"

# Row 45 - CoinLion
$ws.Range("D45").Value = "'0.00005431"
$ws.Range("E45").Value = "'0.27%"

# Row 46 - Kangarootoken (volume only)
$ws.Range("E46").Value = "'0.16%"

# Row 47 - CoinbaseStockToken (volume only)
$ws.Range("E47").Value = "'3.08%"

# Row 48 - BOLO
$ws.Range("D48").Value = "'0.1048"
$ws.Range("E48").Value = "'4,801.48%"

# Row 49 - CryptobidCoin
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.16%"

# Row 50 - SpecialPowerGold
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'0.16%"
